# "Đơn phụ phẫu 1" is worksheet #2 (1-based) in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# 1) Ngày thực hiện của dòng 2 đổi từ 08-02-2024 -> 08-03-2024 (giữ dạng text).
$ws.Range("C2").Value = "'08-03-2024"

# 2) Chèn 1 dòng dữ liệu mới phía trên dòng "Tổng" (đẩy dòng Tổng từ hàng 5 xuống hàng 6).
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "HD-LUXURY"
$ws.Range("B5").Value = 636
$ws.Range("C5").Value = "'08-09-2024"
$ws.Range("D5").Value = "SÓC TRĂNG"
$ws.Range("E5").Value = "thạch thị siêu"
$ws.Range("F5").Value = "CTV"
$ws.Range("G5").Value = "cấy mỡ mặt "
$ws.Range("H5").Value = "Kha Như Huỳnh "

# 3) Cập nhật dòng "Tổng" (giờ ở hàng 6): số đơn tăng từ 3 -> 4, tổng công phụ phẫu 1 giữ nguyên 200000.
$ws.Range("B6").Value = 4
$ws.Range("I6").Value = 200000

Write-Host "edit applied"
